# Insert a new weekly price row for "Apio" (Terminal La Palmera de La Serena)
# right before the existing row 619, shifting all subsequent rows down by one
# (old row 619 becomes 620, ..., old row 670 becomes 671).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 619:670 down to 620:671, leaving row 619 empty and ready to fill.
$ws.Rows(619).Insert()

# Populate the newly inserted row 619 with the new weekly record.
$ws.Cells.Item(619, 1).Value  = 8
$ws.Cells.Item(619, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(619, 3).Value  = "Coquimbo"
$ws.Cells.Item(619, 4).Value  = 45106
$ws.Cells.Item(619, 5).Value  = 4
$ws.Cells.Item(619, 6).Value  = 100112017
$ws.Cells.Item(619, 7).Value  = "Apio"
$ws.Cells.Item(619, 8).Value  = "Americana (o)"
$ws.Cells.Item(619, 9).Value  = "Primera"
$ws.Cells.Item(619, 10).Value = 1300
$ws.Cells.Item(619, 11).Value = 7500
$ws.Cells.Item(619, 12).Value = 8000
$ws.Cells.Item(619, 13).Value = 7750
$ws.Cells.Item(619, 14).Value = "`$/docena de matas"
$ws.Cells.Item(619, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(619, 16).Value = 1292
$ws.Cells.Item(619, 17).Value = 6
$ws.Cells.Item(619, 18).Value = "Hortaliza"
